# Change the "View" field header to "Cache" and reset the default
# value of the boolean column (F) for every data row to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Header rename: F1 "View" -> "Cache"
$ws.Range("F1").Value = "Cache"

# Determine the last used row in column A and set every F-column
# cell (rows 2..last) to FALSE.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
}
